# BOM.xlsx update: "Add files via upload" / "Updated BOM to include drill"
#
# The TOOL LIST section gets a new line item ("Drill with various size
# meshes") inserted right after "Screwdrivers" (and before "Soldering iron
# and tin"). That pushes every following row down by one. The footer lines
# at the bottom of the sheet are updated too: "last update: 2023-02-21"
# becomes "last update: 2023-06-22" and "Rev 1.01" becomes "Rev 1.02".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row right before the current "Soldering iron and tin" row
# (row 73), shifting it and everything below down by one row.
$ws.Rows.Item(73).Insert()

# Fill in the new tool-list entry with the same look as its neighbours.
$ws.Range("A73").Value = "Drill with various size meshes"
$ws.Range("A73").Font.Name = $ws.Range("A72").Font.Name
$ws.Range("A73").Font.Size = $ws.Range("A72").Font.Size
$ws.Range("A73").Font.Bold = $ws.Range("A72").Font.Bold
$ws.Range("A73:G73").Interior.ColorIndex = $ws.Range("A72:G72").Interior.ColorIndex

# Update the revision footer (now two rows further down because of the
# row we just inserted).
$ws.Range("A88").Value = "last update: 2023-06-22"
$ws.Range("A89").Value = "Rev 1.02"

# Match the author's on-save selection/viewport (cursor ended up on the
# updated revision cell after scrolling down to it).
$ws.Activate()
$ws.Range("A89").Select()
